$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''Vans Old Skool'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'
$ws.Range("E3").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''Puma Suede Classic'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'
$ws.Range("E4").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''Converse Chuck 70'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'
$ws.Range("E5").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''Nike Air Force 1'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'
$ws.Range("E6").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''New Balance 550'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'
$ws.Range("E7").Value = 'Lỗi: Expected condition failed: waiting for element to be clickable: By.xpath: //div[contains(@class,''product-card'')][.//a[contains(normalize-space(.),''Adidas Superstar'')]]//button[contains(@class,''btn-dark'') or contains(normalize-space(.),''MUA'')] (tried for 15 second(s) with 500 milliseconds interval)
Build info: version: ''4.14.1'', revision: ''03f8ede370''
System info: os.name: ''Windows 11'', os.arch: ''amd64'', os.version: ''10.0'', java.version: ''17.0.10''
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, browserName: chrome, browserVersion: 145.0.7632.111, chrome: {chromedriverVersion: 145.0.7632.117 (de4e5d6e13b..., userDataDir: C:\Users\tanng\AppData\Loca...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:63620}, goog:processID: 16052, networkConnectionEnabled: false, pageLoadStrategy: normal, platformName: windows, proxy: Proxy(), se:cdp: ws://localhost:63620/devtoo..., se:cdpVersion: 145.0.7632.111, setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}
Session ID: 0015102cfb92e088fda97195c3648cd3'

$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()
$ws.Rows.Item(5).EntireRow.AutoFit()
$ws.Rows.Item(6).EntireRow.AutoFit()
$ws.Rows.Item(7).EntireRow.AutoFit()

$wb.Save()
